# Aleksin test Dashboard General-Dashboard-Page_overview
#
# Adds a new test-data row ("General-Dashboard-Page_overview_[WEB]" / "C70774")
# to the bottom of the table on Sheet1, mirroring the formatting of the
# previous row, and grows the AutoFilter range by one row (matching the
# pre-existing pattern where the filter range lags one row behind the
# freshly-typed row of data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Resize the AutoFilter from A1:F14 to A1:F15 --------------------
# Do this BEFORE inserting the new row, while the worksheet's used range
# still ends at row 15, so the engine's auto-extend-to-used-range behavior
# lands exactly on row 15 (matching the target ref="A1:F15").
$ws.AutoFilterMode = $false
$ws.Range("A1:F15").AutoFilter()

# Keep the hidden _xlnm._FilterDatabase defined name synced with the
# resized AutoFilter range (Excel normally keeps these in lock-step).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$F`$15"
    }
}

# --- 2. Insert the new data row (row 16) --------------------------------
# Copy row 15 (which carries the shared C/D/E/F style+values for this
# table) and insert it as a new row 16, which shifts nothing below it
# (there is nothing below) while preserving formatting exactly.
$ws.Rows("15").Copy()
$ws.Rows("16").Insert()

# Overwrite the two test-specific columns with the new row's content.
$ws.Range("A16").Value = "General-Dashboard-Page_overview_[WEB]"
$ws.Range("B16").Value = "C70774"

# --- 3. Restore the author's final selection ----------------------------
$ws.Range("A19").Select()
